# Insert a new row at row 36 of the "2024" sheet, shifting the existing
# rows 36-111 down to 37-112 (dimension grows from A1:Y111 to A1:Y112).
# The newly inserted row 36 records a fresh "Others" entry that duplicates
# the September details text of the (old) row 36 / (new) row 37 entry
# ("bal axisbank") but stamps it with the timestamp that immediately
# precedes it chronologically (the same timestamp already present on the
# row above, row 35: "2024-09-09 11:38:16").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Shift rows 36:111 down to 37:112 by inserting a blank row at 36.
$ws.Rows("36").Insert()

# Populate the new row's September columns (R = September_Details,
# S = September_Date).
$ws.Range("R36").Value = "bal axisbank"
$ws.Range("S36").Value = "2024-09-09 11:38:16"
